$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.424.75"
$ws.Range("E2").Value = "  -0.62%  "
$ws.Range("D3").Value = "3.436.92"
$ws.Range("E3").Value = "  -1.09%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'407.39"
$ws.Range("E5").Value = "  -0.28%  "
$ws.Range("D6").Value = "'134.59"
$ws.Range("E6").Value = "  +2.79%  "
$ws.Range("E7").Value = "  -0.85%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").Value = "'0.685"
$ws.Range("E9").Value = "  -0.56%  "
$ws.Range("E10").Value = "  -3.84%  "
$ws.Range("D11").Value = "'42.39"
$ws.Range("E11").Value = "  -0.52%  "
$ws.Range("E12").Value = "  -0.63%  "
$ws.Range("D13").Value = "'8.45"
$ws.Range("E13").Value = "  -2.77%  "
$ws.Range("D14").Value = "'19.94"
$ws.Range("E14").Value = "  -0.61%  "
$ws.Range("D15").Value = "3.444.45"
$ws.Range("E15").Value = "  -2.31%  "
$ws.Range("D16").Value = "62.320.94"
$ws.Range("E16").Value = "  -0.60%  "
$ws.Range("D17").Value = "'11.46"
$ws.Range("E17").Value = "  +5.00%  "
$ws.Range("E18").Value = "  -2.51%  "
$ws.Range("D19").Value = "'0.0000131"
$ws.Range("E19").Value = "  -2.76%  "
$ws.Range("D20").Value = "'3.20"
$ws.Range("E20").Value = "  -4.66%  "
$ws.Range("D21").Value = "'84.18"
$ws.Range("E21").Value = "  +2.01%  "
$ws.Range("D22").Value = "'314.72"
$ws.Range("E22").Value = "  +1.98%  "
$ws.Range("D23").Value = "'12.97"
$ws.Range("E23").Value = "  -1.14%  "
$ws.Range("D24").Value = "'3.16"
$ws.Range("E24").Value = "  -0.39%  "
$ws.Range("D25").Value = "'4.75"
$ws.Range("E25").Value = "  +8.39%  "
$ws.Range("E26").Value = "  -1.84%  "
$ws.Range("D27").Value = "'8.25"
$ws.Range("E27").Value = "  -0.91%  "
$ws.Range("B28").Value = "RenderToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D28").Value = "'7.62"
$ws.Range("E28").Value = "  -1.92%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "'2.73"
$ws.Range("E29").Value = "  +2.27%  "
$ws.Range("D30").Value = "'0.173"
$ws.Range("E30").Value = "  -4.20%  "
$ws.Range("E31").Value = "  -3.53%  "
$ws.Range("B32").Value = "Dai"
$ws.Range("C32").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D32").Value = "'1.00"
$ws.Range("E32").Value = "  -0.04%  "
$ws.Range("B33").Value = "InjectiveProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D33").Value = "'42.23"
$ws.Range("E33").Value = "  -2.24%  "
$ws.Range("D34").Value = "'11.39"
$ws.Range("E34").Value = "  -4.24%  "
$ws.Range("D35").Value = "'0.0485"
$ws.Range("E35").Value = "  -1.21%  "
$ws.Range("D36").Value = "'51.50"
$ws.Range("E36").Value = "  -2.12%  "
$ws.Range("D37").Value = "'0.997"
$ws.Range("E37").Value = "  +0.03%  "
$ws.Range("D38").Value = "'3.42"
$ws.Range("E38").Value = "  -4.20%  "
$ws.Range("D39").Value = "'2.96"
$ws.Range("E39").Value = "  -1.18%  "
$ws.Range("D40").Value = "'0.317"
$ws.Range("E40").Value = "  +10.89%  "
$ws.Range("D41").Value = "'137.91"
$ws.Range("E41").Value = "  -0.06%  "
$ws.Range("D42").Value = "'1.99"
$ws.Range("E43").Value = "  -0.27%  "
$ws.Range("D44").Value = "'4.03"
$ws.Range("E44").Value = "  +1.67%  "
$ws.Range("D45").Value = "'16.85"
$ws.Range("E45").Value = "  -3.81%  "
$ws.Range("E46").Value = "  -1.02%  "
$ws.Range("D47").Value = "'21.52"
$ws.Range("E47").Value = "  -3.75%  "
$ws.Range("D48").Value = "2.129.09"
$ws.Range("E48").Value = "  -3.53%  "
$ws.Range("E49").Value = "  -4.51%  "
$ws.Range("D50").Value = "'1.92"
$ws.Range("E50").Value = "  +3.35%  "
$ws.Range("D51").Value = "'1.73"
$ws.Range("E51").Value = "  +21.27%  "
